$d = $word.ActiveDocument
$blue = 0x00B44600   # BGR encoding of RGB 0046B4, used by Word's Font.Color

# ---------------------------------------------------------------------------
# 1. Split the "CATALOG NO / LOT NO" paragraph into two separate paragraphs,
#    right after the catalog number ("IMSKLK1KT"), so "   LOT NO: 20250424"
#    becomes its own paragraph.
# ---------------------------------------------------------------------------
$splitRng = $d.Content
$splitRng.Find.Execute("IMSKLK1KT", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$splitRng.Collapse(0)
$splitRng.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 2. Re-locate the two paragraphs by their known text content.
# ---------------------------------------------------------------------------
$catalogPara = $null
$lotPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "CATALOG NO:*") {
        $catalogPara = $p
    } elseif ($t -like "*LOT NO:*") {
        $lotPara = $p
    }
}

# ---------------------------------------------------------------------------
# 3. Apply the Heading2 paragraph style to both paragraphs (this also drops
#    the now-redundant direct <w:b/> on the label runs, since Heading2
#    already renders bold).
# ---------------------------------------------------------------------------
$catalogPara.Range.Style = "Heading2"
$lotPara.Range.Style = "Heading2"

# ---------------------------------------------------------------------------
# 4. Color the "CATALOG NO: " and "LOT NO: " label runs with the new blue
#    accent color.
# ---------------------------------------------------------------------------
$labelRng1 = $d.Content
$labelRng1.Find.Execute("CATALOG NO: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$labelRng1.Font.Color = $blue

$labelRng2 = $d.Content
$labelRng2.Find.Execute("LOT NO: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$labelRng2.Font.Color = $blue

# ---------------------------------------------------------------------------
# 5. Remove the leading spaces before "LOT NO: " so the run text becomes
#    exactly "LOT NO: " (it used to be "   LOT NO: ").
# ---------------------------------------------------------------------------
$spacesRng = $d.Content
$spacesRng.Find.Execute("   LOT NO: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$delRng = $d.Range($spacesRng.Start, $spacesRng.Start + 3)
$delRng.Delete()

# ---------------------------------------------------------------------------
# 6. Set the explicit alignment: left for the catalog line, right for the
#    lot line.
# ---------------------------------------------------------------------------
$catalogPara.Range.ParagraphFormat.Alignment = 0   # wdAlignParagraphLeft
$lotPara.Range.ParagraphFormat.Alignment = 2        # wdAlignParagraphRight

# ---------------------------------------------------------------------------
# 7. Update the Heading2 style definition itself: its default run color
#    moves from dark navy (00008B) to the new accent blue (0046B4).
# ---------------------------------------------------------------------------
$heading2Style = $d.Styles.Item("Heading 2")
$heading2Style.Font.Color = $blue
